# Fix the "Dewst_template" typo -> "Dest_template" on the four section-divider
# slides (slide 2, 4, 5 and 7) and bump the trailing slide counter that
# follows the en-dash by one (Slide2->Slide2(unchanged text itself) /
# Slide2->Slide3 / Slide3->Slide4 / Slide4->Slide5), splitting the run that
# holds " <en-dash> " so the leading space becomes its own run, matching how
# PowerPoint re-segments runs when you edit text in place.

$p = $ppt.ActivePresentation

function Fix-TitleRun($Shape, $NewSuffix) {
    $tr = $Shape.TextFrame.TextRange

    # 1) "Dewst_template" -> "Dest_template" (remove the stray 'w').
    #    Re-assigning the whole first run's Text keeps it as a single run
    #    instead of splintering it character by character.
    $run1 = $tr.Runs(1)
    $run1.Text = "Dest_template"

    # 2) The run that used to read " <en-dash> " (or " <en-dash> SlideN")
    #    should be split so the leading space becomes its own run, and the
    #    trailing "SlideN" text becomes/ends up as its own run too.
    #    Re-fetch the TextRange after each structural edit.
    $tr = $Shape.TextFrame.TextRange
    $tr.Characters(14, 1).Text = " "

    # 3) Replace the trailing "SlideN" text with the new slide number.
    $tr = $Shape.TextFrame.TextRange
    $fullLen = $tr.Length
    $tail = $tr.Characters(17, $fullLen - 16)
    $tail.Text = $NewSuffix
}

# Slide 2: "Dewst_template - Slide2" -> "Dest_template - Slide2"
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(3)
Fix-TitleRun $sh2 "Slide2"

# Slide 4: "Dewst_template - Slide2" -> "Dest_template - Slide3"
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(5)
Fix-TitleRun $sh4 "Slide3"

# Slide 5: "Dewst_template - Slide3" -> "Dest_template - Slide4"
$s5 = $p.Slides.Item(5)
$sh5 = $s5.Shapes.Item(5)
Fix-TitleRun $sh5 "Slide4"

# Slide 7: "Dewst_template - Slide4" -> "Dest_template - Slide5"
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(5)
Fix-TitleRun $sh7 "Slide5"
